{"js": "// The document is a single-column table where each row's one cell holds a\n// single stat value (a couple of summary rows hold multiple tab-separated\n// values inside one run). This edit updates a handful of those values and\n// collapses the three trailing multi-value rows down to a single value each\n// (matching the \"Fixed README.md stats\" commit).\n\nconst table = context.document.body.tables.getFirstOrNullObject();\ntable.load(\"items\");\nawait context.sync();\n\n// index -> new cell text (0-based row index within the single table)\nconst updates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"164\",\n  5: \"0.00058\",\n  6: \"0.00017\",\n  7: \"0.00006\",\n  8: \"0.00022\",\n  9: \"0.00030\",\n  10: \"0.00043\",\n  11: \"0.02850\",\n  // These three rows currently hold many tab-separated <w:t> runs; setting\n  // `.value` replaces the whole cell range with a single run, matching the\n  // diff which drops the tab-separated values in favor of one short value.\n  43: \"99.99\",\n  44: \"0.03\",\n  45: \"419\",\n};\n\nfor (const [rowIndex, text] of Object.entries(updates)) {\n  const cell = table.getCell(Number(rowIndex), 0);\n  cell.value = text;\n}\n\nawait context.sync();\n", "ps1": "# The document is a single-column table where each row's one cell holds a\n# single stat value (a couple of summary rows hold multiple tab-separated\n# values inside one run). This edit updates a handful of those values and\n# collapses the three trailing multi-value rows down to a single value each\n# (matching the \"Fixed README.md stats\" commit).\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# Word's Cell()/Rows collections are 1-based, so row N here is the (N-1)-th\n# 0-based table row discussed in the diff.\n$updates = [ordered]@{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"164\"\n    6  = \"0.00058\"\n    7  = \"0.00017\"\n    8  = \"0.00006\"\n    9  = \"0.00022\"\n    10 = \"0.00030\"\n    11 = \"0.00043\"\n    12 = \"0.02850\"\n    # These three rows currently hold many tab-separated runs in one cell;\n    # writing Range.Text replaces the whole cell content with a single run,\n    # matching the diff which drops the tab-separated values for one value.\n    44 = \"99.99\"\n    45 = \"0.03\"\n    46 = \"419\"\n}\n\nforeach ($rowIndex in $updates.Keys) {\n    $cell = $table.Cell($rowIndex, 1)\n    $cell.Range.Text = $updates[$rowIndex]\n}\n"}
